# Auto-generated Excel COM-interop script to apply weekly crime-data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (title volume/number + reporting week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Crime statistics table updates ---

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -36.363636363636
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = -16.666666666666
$ws.Range("L16").Value = -28.571428571428
$ws.Range("M16").Value = -41.176470588235
$ws.Range("N16").Value = -81.481481481481

# Row 17
$ws.Range("C17").Value = "'0"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = -55.555555555555
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = -54.545454545454
$ws.Range("L17").Value = -50
$ws.Range("N17").Value = -77.272727272727

# Row 18
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 16
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 16
$ws.Range("L18").Value = -40.74074074074
$ws.Range("M18").Value = -11.111111111111
$ws.Range("N18").Value = -86.776859504132

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -19.230769230769
$ws.Range("I19").Value = 50
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = -15.254237288135
$ws.Range("L19").Value = -10.714285714285
$ws.Range("M19").Value = 108.333333333333
$ws.Range("N19").Value = 25

# Row 20
$ws.Range("F20").Value = 1
$ws.Range("H20").Value = -80
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = -40
$ws.Range("M20").Value = -57.142857142857
$ws.Range("N20").Value = -96.103896103896

# Row 21
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -10
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -23.076923076923
$ws.Range("I21").Value = 84
$ws.Range("J21").Value = 103
$ws.Range("K21").Value = -18.446601941747
$ws.Range("L21").Value = -26.95652173913
$ws.Range("M21").Value = 18.309859154929
$ws.Range("N21").Value = -73.417721518987

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 400

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("I23").Value = 3
$ws.Range("K23").Value = 200
$ws.Range("L23").Value = 50

# Row 24
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 2.898550724637
$ws.Range("I24").Value = 79
$ws.Range("J24").Value = 82
$ws.Range("K24").Value = -3.658536585365
$ws.Range("L24").Value = 9.722222222222
$ws.Range("M24").Value = 97.5

# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 114.285714285714
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 38.235294117647
$ws.Range("I25").Value = 56
$ws.Range("J25").Value = 38
$ws.Range("K25").Value = 47.368421052631
$ws.Range("L25").Value = 40

# Row 26
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -57.142857142857
$ws.Range("F26").Value = 14
$ws.Range("H26").Value = -39.130434782608
$ws.Range("I26").Value = 17
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -37.037037037037
$ws.Range("L26").Value = -26.086956521739
$ws.Range("M26").Value = 13.333333333333

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 4
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = 300
$ws.Range("L28").Value = -42.857142857142

# Row 31
$ws.Range("L31").Value = -100
